$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Insert a new row above row 31 by copying row 31 (this duplicates
# the existing "danholland" row's formatting/values, pushing it and
# everything below down by one), then overwrite the new row's
# UserName/Password with the new Oleg_Babak account.
$ws.Rows.Item(31).Copy()
$ws.Rows.Item(31).Insert(-4161)

$ws.Cells.Item(31, 1).Value = "Oleg_Babak"
$ws.Cells.Item(31, 2).Value = "Password1!"

# Restore the view state to match where the author was working when
# they saved the workbook.
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Range("G30").Select()
